$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = "-"
$ws.Range("C3").Value = "-"
$ws.Range("F3").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("C4").Value = "-"
$ws.Range("F4").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("F6").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("F7").Value = "['MCT-3A-CAM', -, -, -]"
$ws.Range("C8").Value = "-"
